$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.120310664176941
$ws.Range("B1").Value = 2.2652268409729
$ws.Range("C1").Value = 10.84642601013184
$ws.Range("D1").Value = 1.801810503005981
$ws.Range("E1").Value = 1.288565158843994
